$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("Tipo") to make room for "MAE".
# This shifts the existing "Tipo" column (D -> E) along with its data.
$ws.Range("D1").EntireColumn.Insert()

# Fill in the new "MAE" column header and value.
$ws.Range("D1").Value = "MAE"
$ws.Range("D2").Value = 0.6472536726993832

# Copy the header style (bold, bordered, centered) from an existing header
# cell so the new header matches the rest of the header row.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
